$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row2 = $ws.Range("A2:R2")
$row2.NumberFormat = "@"

$ws.Range("A2").Value = "Practice"
$ws.Range("B2").Value = "23"
$ws.Range("C2").Value = "24255"
$ws.Range("D2").Value = "88"
$ws.Range("E2").Value = "1"
$ws.Range("F2").Value = "15"
$ws.Range("G2").Value = "Specimen"
$ws.Range("H2").Value = "150"
$ws.Range("I2").Value = "15"
$ws.Range("J2").Value = "2"
$ws.Range("K2").Value = "2"
$ws.Range("L2").Value = "Tier 1"
$ws.Range("M2").Value = "34"
$ws.Range("N2").Value = "15"
$ws.Range("O2").Value = "5"
$ws.Range("P2").Value = "No"
$ws.Range("Q2").Value = "No"
$ws.Range("R2").Value = "This team is awesome. A lot of potential. "

$row2.ClearFormats()
